$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.01955
$ws.Range("E2").Value = 0.03856
$ws.Range("F2").Value = 0.062
$ws.Range("G2").Value = 0.0846659311419022
$ws.Range("H2").Value = 0.0846659311419022
$ws.Range("I2").Value = 0.06338160586054123
$ws.Range("J2").Value = 0.05597745668061119
$ws.Range("K2").Value = 2649.1
$ws.Range("L2").Value = 0.04298114676964014
$ws.Range("M2").Value = 2761.7
$ws.Range("N2").Value = 0.1235516226300296
$ws.Range("O2").Value = 1.04250500169869
$ws.Range("P2").Value = 497.4
$ws.Range("Q2").Value = 0.02225244490573804
$ws.Range("R2").Value = 0.1877618813936808
$ws.Range("S2").Value = 2264.3
$ws.Range("T2").Value = 0.819893543831698
$ws.Range("U2").Value = 26475.3
$ws.Range("V2").Value = 1.184439394074962
$ws.Range("W2").Value = 0.04084668672067628
$ws.Range("X2").Value = 0.09766601826728588
$ws.Range("Y2").Value = -0.0568193315466096
$ws.Range("Z2").Value = 0.9800817549248489
$ws.Range("AA2").Value = 0.05137411929916241
$ws.Range("AB2").Value = 0.05310781966187775
$ws.Range("AC2").Value = -0.001733700362715338
$ws.Range("AD2").Value = 25555.9
$ws.Range("AE2").Value = 222.6905219570081
$ws.Range("AF2").Value = 25778.59052195701
$ws.Range("AG2").Value = -696.7094780429907
$ws.Range("AH2").Value = 0.5355901286131091
$ws.Range("AI2").Value = 0.2672542128011587
$ws.Range("AJ2").Value = -0.03217182305833111
$ws.Range("AK2").Value = -0.009955568116931567
$ws.Range("AL2").Value = 1174.9
$ws.Range("AM2").Value = 1174.9
$ws.Range("AN2").Value = 6.289134982158238
$ws.Range("AO2").Value = 3.316026895906034
$ws.Range("AP2").Value = -0.1714555132380929
$ws.Range("AQ2").Value = 3.316026895906034

# Row 3
$ws.Range("D3").Value = 0.076
$ws.Range("E3").Value = 0.00042
$ws.Range("F3").Value = 0.18
$ws.Range("G3").Value = 0.1510326848003746
$ws.Range("H3").Value = 0.1510326848003746
$ws.Range("I3").Value = 0.1221166098581212
$ws.Range("J3").Value = 0.09592405967588871
$ws.Range("K3").Value = 1610.4
$ws.Range("L3").Value = 0.07251081763781514
$ws.Range("M3").Value = 1140
$ws.Range("N3").Value = 0.07941760423560557
$ws.Range("O3").Value = 0.7078986587183308
$ws.Range("P3").Value = 249.8
$ws.Range("Q3").Value = 0.01740220836671427
$ws.Range("R3").Value = 0.1551167411823149
$ws.Range("S3").Value = 890.2
$ws.Range("T3").Value = 0.7808771929824562
$ws.Range("U3").Value = 13941
$ws.Range("V3").Value = 0.9711937023233133
$ws.Range("W3").Value = 0.0440206873175372
$ws.Range("X3").Value = 0.08553919197195863
$ws.Range("Y3").Value = -0.04151850465442142
$ws.Range("Z3").Value = 0.6138501934770589
$ws.Range("AA3").Value = 0.05888300259114923
$ws.Range("AB3").Value = 0.05366944207214916
$ws.Range("AC3").Value = 0.005213560519000071
$ws.Range("AD3").Value = 12906.8
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 12906.8
$ws.Range("AG3").Value = -1034.200000000001
$ws.Range("AH3").Value = 0.473447707922953
$ws.Range("AI3").Value = 0.2359582811543067
$ws.Range("AJ3").Value = -0.07764089397385951
$ws.Range("AK3").Value = -0.02537384839972032
$ws.Range("AL3").Value = 608.8
$ws.Range("AM3").Value = 608.8
$ws.Range("AN3").Value = 4.624601383066394
$ws.Range("AO3").Value = 4.454829172141919
$ws.Range("AP3").Value = -0.3705614676269307
$ws.Range("AQ3").Value = 4.454829172141919

# Row 4
$ws.Range("D4").Value = -0.0369
$ws.Range("E4").Value = 0.0767
$ws.Range("F4").Value = -0.05599999999999999
$ws.Range("G4").Value = 0.04727976481868058
$ws.Range("H4").Value = 0.04727976481868058
$ws.Range("I4").Value = 0.0302946081184378
$ws.Range("J4").Value = 0.02971449252103728
$ws.Range("K4").Value = 1038.7
$ws.Range("L4").Value = 0.02634629383968
$ws.Range("M4").Value = 1621.7
$ws.Range("N4").Value = 0.2027606556557182
$ws.Range("O4").Value = 1.561278521228458
$ws.Range("P4").Value = 247.6
$ws.Range("Q4").Value = 0.03095735237118815
$ws.Range("R4").Value = 0.2383748916915375
$ws.Range("S4").Value = 1374.1
$ws.Range("T4").Value = 0.8473207128322131
$ws.Range("U4").Value = 12534.3
$ws.Range("V4").Value = 1.567159700428852
$ws.Range("W4").Value = 0.03767268612381536
$ws.Range("X4").Value = 0.1097928445626131
$ws.Range("Y4").Value = -0.07212015843879777
$ws.Range("Z4").Value = 1.476223629803533
$ws.Range("AA4").Value = 0.04386523600717559
$ws.Range("AB4").Value = 0.05254619725160634
$ws.Range("AC4").Value = -0.008680961244430747
$ws.Range("AD4").Value = 12649.1
$ws.Range("AE4").Value = 222.6905219570081
$ws.Range("AF4").Value = 12871.79052195701
$ws.Range("AG4").Value = 337.4905219570101
$ws.Range("AH4").Value = 0.6167636820334407
$ws.Range("AI4").Value = 0.3082495789653111
$ws.Range("AJ4").Value = 0.04048789597665779
$ws.Range("AK4").Value = 0.0115486436012083
$ws.Range("AL4").Value = 566.1
$ws.Range("AM4").Value = 566.1
$ws.Range("AN4").Value = 9.93957252868144
$ws.Range("AO4").Value = 2.091326620738386
$ws.Range("AP4").Value = 0.2651976441592095
$ws.Range("AQ4").Value = 2.091326620738386
